# Africa Cup of Nations and Copa Centroamericana 2017
# Insert a new "group_stage" table at the top of the 2017 sheet, above the
# existing "game" and "game_score" tables (which get pushed down by 8 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2017")

# Push the existing two tables ("game" headers/data starting at row 1, and
# "game_score" headers/data starting at row 18) down by 8 rows, opening up
# rows 1-8 for the new "group_stage" table (row 8 stays blank, mirroring the
# existing blank-row separator pattern used between tables on this sheet).
$ws.Rows("1:8").Insert()

# --- Row 1: group_stage header row ---
$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "tournament"
$ws.Cells.Item(1,3).Value = "group_code"
$ws.Cells.Item(1,4).Value = "squad"
$ws.Cells.Item(1,1).HorizontalAlignment = -4152
$ws.Cells.Item(1,7).Formula = '="insert into group_stage (id, tournament, group_code, squad) values (" & A1 & ", " & B1 & ", ''" & C1 & "'', " & D1 &  ");"'

# --- Row 2: first group_stage data row ---
$ws.Cells.Item(2,1).Formula = "='2014'!A8+1"
$ws.Cells.Item(2,2).Value = 2017
$ws.Cells.Item(2,3).Value = "A"
$ws.Cells.Item(2,4).Value = 505
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,7).Formula = '="insert into group_stage (id, tournament, group_code, squad) values (" & A2 & ", " & B2 & ", ''" & C2 & "'', " & D2 &  ");"'

# --- Row 3 ---
$ws.Cells.Item(3,1).Formula = "=A2+1"
$ws.Cells.Item(3,2).Formula = "=B2"
$ws.Cells.Item(3,3).Value = "A"
$ws.Cells.Item(3,4).Value = 503
$ws.Cells.Item(3,7).Formula = '="insert into group_stage (id, tournament, group_code, squad) values (" & A3 & ", " & B3 & ", ''" & C3 & "'', " & D3 &  ");"'

# --- Row 4 ---
$ws.Cells.Item(4,1).Formula = "=A3+1"
$ws.Cells.Item(4,2).Formula = "=B3"
$ws.Cells.Item(4,3).Value = "A"
$ws.Cells.Item(4,4).Value = 504
$ws.Cells.Item(4,7).Formula = '="insert into group_stage (id, tournament, group_code, squad) values (" & A4 & ", " & B4 & ", ''" & C4 & "'', " & D4 &  ");"'

# --- Row 5 ---
$ws.Cells.Item(5,1).Formula = "=A4+1"
$ws.Cells.Item(5,2).Formula = "=B4"
$ws.Cells.Item(5,3).Value = "A"
$ws.Cells.Item(5,4).Value = 501
$ws.Cells.Item(5,7).Formula = '="insert into group_stage (id, tournament, group_code, squad) values (" & A5 & ", " & B5 & ", ''" & C5 & "'', " & D5 &  ");"'

# --- Row 6 ---
$ws.Cells.Item(6,1).Formula = "=A5+1"
$ws.Cells.Item(6,2).Formula = "=B5"
$ws.Cells.Item(6,3).Value = "A"
$ws.Cells.Item(6,4).Value = 506
$ws.Cells.Item(6,7).Formula = '="insert into group_stage (id, tournament, group_code, squad) values (" & A6 & ", " & B6 & ", ''" & C6 & "'', " & D6 &  ");"'

# --- Row 7 ---
$ws.Cells.Item(7,1).Formula = "=A6+1"
$ws.Cells.Item(7,2).Formula = "=B6"
$ws.Cells.Item(7,3).Value = "A"
$ws.Cells.Item(7,4).Value = 507
$ws.Cells.Item(7,7).Formula = '="insert into group_stage (id, tournament, group_code, squad) values (" & A7 & ", " & B7 & ", ''" & C7 & "'', " & D7 &  ");"'

# Row 8 is left blank, matching the blank separator row that already
# existed between the "game" and "game_score" tables.
